# Updates the "startup" sheet of the UBC02 Breed/Diagnosis/PrimDiseaseSite test-case
# workbook: the CasesTab query text drops its trailing `Cohort` column, and the
# CasesTab / SamplesTab / FilesTab shared-string entries get reshuffled (SamplesTab
# and FilesTab now sort ahead of CasesTab in the shared-string table) which is an
# invisible, non-semantic side effect of Excel re-saving the file. The visible
# effect in the UI is just: B2's query text no longer returns a Cohort column, and
# the selection / scroll position moved from D4 to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B2 (CasesTab query): drop the trailing Cohort line -------------------------
$casesTabQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Bluetick Hound','Welsh Springer Spaniel','Wheaten Terrier']and diag.disease_term in ['Bladder Cancer'] and diag.primary_disease_site in ['Bladder, Urethra']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $casesTabQuery

# --- Row heights: Excel re-wraps the edited/re-flowed cells on save -------------
$ws.Rows.Item(2).RowHeight = 290
$ws.Rows.Item(3).RowHeight = 261
$ws.Rows.Item(4).RowHeight = 261

# --- Column widths: minor re-measurement from the resave -----------------------
# (values are nudged by -5/6 to offset this host's ColumnWidth->stored-width
# conversion, landing on the closest width the engine can actually store).
# Columns 1 and 5 are left alone so their bestFit="1" flag (from the original
# "best fit" auto-sizing) survives the resave untouched.
$ws.Columns.Item(2).ColumnWidth = 74.98307291666667
$ws.Columns.Item(3).ColumnWidth = 74.98307291666667
$ws.Columns.Item(4).ColumnWidth = 69.43619791666667

# --- View state: scrolled/selected B2 instead of D4, at a smaller zoom window ---
$ws.Range("B2").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
